# The edit reorders the "Artfynd" species-observation rows (2,3,4 and
# 7,8,9,10 - rows 5 and 6 are left untouched): each destination row ends
# up holding the field values that another source row held beforehand.
# Concretely: new2<-old3, new3<-old10, new10<-old2 (a 3-cycle), and
# new4<->old8, new7<->old9 (two swaps). Rather than moving cells with
# Cut/Insert (which would also drag row formatting/height along), we
# write the literal target values straight into each destination cell -
# that reproduces the same end state the author's commit shows.
#
# Notes on COM quirks in this runtime:
#  - Plain `.Value = "12"` on a digit-only string gets auto-coerced to a
#    Number by the COM layer. To keep it stored as Text (as the source
#    file does for columns like "Antal"), we prefix with a text-quote
#    (`'`) and then reset `.Style = "Normal"` so the quote-prefix flag
#    doesn't leave a stray cell style behind.
#  - A handful of cells must end up present-but-empty (empty string,
#    Text-typed) rather than absent. `ClearContents()` drops the cell
#    entirely (becomes a blank/Number-typed cell on save), so for those
#    we instead write a lone text-quote (`'`) which Excel stores as an
#    empty string while keeping the cell's Text type, then reset the
#    style the same way.
#  - Cells that should become entirely absent (no more value at all) use
#    `ClearContents()`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112128588
$ws.Range("B2").Value = 5113
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 100526
$ws.Range("F2").Value = "Bronshjon"
$ws.Range("G2").Value = "Callidium coriaceum"
$ws.Range("H2").Value = "Paykull, 1800"
$ws.Range("I2").Value = "'1"
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = "'"
$ws.Range("J2").Style = "Normal"
$ws.Range("M2").Value = "färska gnagspår"
$ws.Range("Q2").Value = 655234
$ws.Range("R2").Value = 6675166
$ws.Range("Z2").Value = "11:08"
$ws.Range("AB2").Value = "11:08"
$ws.Range("AC2").ClearContents()

# Row 3
$ws.Range("A3").Value = 112128602
$ws.Range("B3").Value = 56575
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 103021
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("J3").ClearContents()
$ws.Range("M3").Value = "permanent revir"
$ws.Range("Q3").Value = 655214
$ws.Range("R3").Value = 6675119
$ws.Range("Z3").Value = "10:09"
$ws.Range("AB3").Value = "10:09"

# Row 4
$ws.Range("A4").Value = 112129079
$ws.Range("B4").Value = 98980
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."
$ws.Range("I4").Value = "'"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = "'"
$ws.Range("J4").Style = "Normal"
$ws.Range("Q4").Value = 655188
$ws.Range("R4").Value = 6675131
$ws.Range("Z4").Value = "10:28"
$ws.Range("AB4").Value = "10:28"
$ws.Range("AC4").ClearContents()

# Row 7
$ws.Range("A7").Value = 112129069
$ws.Range("I7").Value = "'2"
$ws.Range("I7").Style = "Normal"
$ws.Range("Q7").Value = 655168
$ws.Range("R7").Value = 6675142
$ws.Range("Z7").Value = "10:37"
$ws.Range("AB7").Value = "10:38"

# Row 8
$ws.Range("A8").Value = 112129067
$ws.Range("B8").Value = 96735
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("I8").Value = "'3"
$ws.Range("I8").Style = "Normal"
$ws.Range("J8").Value = "plantor/tuvor"
$ws.Range("Q8").Value = 655241
$ws.Range("R8").Value = 6675125
$ws.Range("Z8").Value = "10:16"
$ws.Range("AB8").Value = "10:17"
$ws.Range("AC8").Value = "Tuff tillvaro nära hyggeskanten."

# Row 9
$ws.Range("A9").Value = 112129073
$ws.Range("I9").Value = "'6"
$ws.Range("I9").Style = "Normal"
$ws.Range("Q9").Value = 655180
$ws.Range("R9").Value = 6675133
$ws.Range("Z9").Value = "10:43"
$ws.Range("AB9").Value = "10:44"

# Row 10
$ws.Range("A10").Value = 112129065
$ws.Range("B10").Value = 96735
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("I10").Value = "'12"
$ws.Range("I10").Style = "Normal"
$ws.Range("J10").Value = "plantor/tuvor"
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value = 655221
$ws.Range("R10").Value = 6675131
$ws.Range("Z10").Value = "10:12"
$ws.Range("AB10").Value = "10:13"
$ws.Range("AC10").Value = "Djupt nere bland ris och mossa."
